$d = $word.ActiveDocument
$tbl = $d.Tables(1)
$values = @(
  "68-49=",
  "78+21=",
  "57+30=",
  "25+0=",
  "76-39=",
  "61-10=",
  "19+20=",
  "24+28=",
  "64-35=",
  "0+74=",
  "19+25=",
  "55-11=",
  "45+39=",
  "2+19=",
  "11+45=",
  "81-33=",
  "35-5=",
  "41+8=",
  "56-52=",
  "75-9=",
  "51-29=",
  "32+6=",
  "89-76=",
  "88-61=",
  "48+33=",
  "84-32=",
  "70-46=",
  "96-59=",
  "66+29=",
  "22+38=",
  "87-41=",
  "75+8=",
  "98-78=",
  "82-35=",
  "3+81=",
  "90-82=",
  "82-81=",
  "85-33=",
  "88-50=",
  "53+33=",
  "33+22=",
  "23+55=",
  "23+28=",
  "77+5=",
  "51-39=",
  "51+24=",
  "17-9=",
  "43-30=",
  "35+11=",
  "66-56=",
  "0+1=",
  "15+80=",
  "34-3=",
  "58-54=",
  "30+63=",
  "17+18=",
  "32+52=",
  "68-67=",
  "62-15=",
  "82-29=",
  "35-15=",
  "90-31=",
  "10+89=",
  "6+63=",
  "36+31=",
  "1+18=",
  "78+10=",
  "3+45=",
  "34+65=",
  "22+42=",
  "7+14=",
  "53-5=",
  "3+51=",
  "31+20=",
  "58+38=",
  "32+41=",
  "36-36=",
  "14+8=",
  "35-33=",
  "7+57=",
  "39-20=",
  "12+56=",
  "87-60=",
  "34-2=",
  "87-12=",
  "29+47=",
  "21+45=",
  "25+18=",
  "58+11=",
  "44+8=",
  "34+0=",
  "51-13=",
  "23+63=",
  "69+19=",
  "66-6=",
  "37-1=",
  "90-66=",
  "96+3=",
  "50+1=",
  "31+24="
)
$rows = $tbl.Rows.Count
$cols = $tbl.Columns.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
  for ($c = 1; $c -le $cols; $c++) {
    $cell = $tbl.Cell($r, $c)
    $rng = $cell.Range
    $rng.MoveEnd(1, -1) | Out-Null
    $rng.Text = $values[$idx]
    $idx = $idx + 1
  }
}
